$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.086.16"
$ws.Range("E2").Value = "  +0.14%  "
$ws.Range("D3").Value = "'2.314.89"
$ws.Range("E3").Value = "  +0.23%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'302.21"
$ws.Range("E5").Value = "  -0.50%  "
$ws.Range("E6").Value = "  -1.54%  "
$ws.Range("D7").Value = "'0.509"
$ws.Range("E7").Value = "  +0.67%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "'0.522"
$ws.Range("E9").Value = "  +2.47%  "
$ws.Range("D10").Value = "'35.90"
$ws.Range("E10").Value = "  +3.05%  "
$ws.Range("D11").Value = "'0.0789"
$ws.Range("E11").Value = "  -0.91%  "
$ws.Range("E12").Value = "  -1.20%  "
$ws.Range("D13").Value = "'17.78"
$ws.Range("E13").Value = "  -0.85%  "
$ws.Range("E14").Value = "  +0.39%  "
$ws.Range("D15").Value = "'2.675.81"
$ws.Range("E15").Value = "  -0.45%  "
$ws.Range("D16").Value = "'2.289.81"
$ws.Range("E16").Value = "  -0.44%  "
$ws.Range("D17").Value = "'0.789"
$ws.Range("E17").Value = "  -3.12%  "
$ws.Range("D18").Value = "'42.997.62"
$ws.Range("E18").Value = "  +0.08%  "
$ws.Range("D19").Value = "'13.29"
$ws.Range("E19").Value = "  +6.50%  "
$ws.Range("D20").Value = "'6.20"
$ws.Range("E20").Value = "  +0.59%  "
$ws.Range("E21").Value = "  +0.36%  "
$ws.Range("D22").Value = "'68.15"
$ws.Range("E22").Value = "  +0.55%  "
$ws.Range("D23").Value = "'240.83"
$ws.Range("E23").Value = "  +1.67%  "
$ws.Range("E24").Value = "  -2.92%  "
$ws.Range("E25").Value = "  -0.74%  "
$ws.Range("E26").Value = "  -0.18%  "
$ws.Range("E27").Value = "  +1.47%  "
$ws.Range("D28").Value = "'169.32"
$ws.Range("E28").Value = "  +0.61%  "
$ws.Range("E29").Value = "  -2.22%  "
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("D31").Value = "'33.34"
$ws.Range("E31").Value = "  -2.15%  "
$ws.Range("D32").Value = "'4.95"
$ws.Range("E32").Value = "  +6.88%  "
$ws.Range("E33").Value = "  +3.30%  "
$ws.Range("E34").Value = "  -0.04%  "
$ws.Range("D35").Value = "'18.31"
$ws.Range("E35").Value = "  +7.04%  "
$ws.Range("E36").Value = "  -0.72%  "
$ws.Range("E37").Value = "  +0.24%  "
$ws.Range("E38").Value = "  +1.43%  "
$ws.Range("D39").Value = "'0.103"
$ws.Range("E39").Value = "  +0.26%  "
$ws.Range("D40").Value = "'2.75"
$ws.Range("E40").Value = "  -2.37%  "
$ws.Range("E41").Value = "  -0.14%  "
$ws.Range("D42").Value = "'1.995.82"
$ws.Range("E42").Value = "  -0.35%  "
$ws.Range("D43").Value = "'0.0290"
$ws.Range("E43").Value = "  +1.08%  "
$ws.Range("D44").Value = "'10.13"
$ws.Range("E44").Value = "  -0.25%  "
$ws.Range("D45").Value = "'2.07"
$ws.Range("E45").Value = "  -10.91%  "
$ws.Range("E46").Value = "  -0.76%  "
$ws.Range("E47").Value = "  -0.39%  "
$ws.Range("D48").Value = "'76.40"
$ws.Range("E48").Value = "  +8.90%  "
$ws.Range("D49").Value = "'54.94"
$ws.Range("E49").Value = "  -1.55%  "
$ws.Range("D50").Value = "'2.541.84"
$ws.Range("E50").Value = "  +0.68%  "
$ws.Range("E51").Value = "  +0.02%  "
